$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.029.03"
$ws.Range("E2").Value = "  -0.88%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.45"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("E4").Value = "  -2.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.53"
$ws.Range("E5").Value = "  -2.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.53"
$ws.Range("E6").Value = "  -3.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.557"
$ws.Range("E7").Value = "  -2.37%  "

# Row 8
$ws.Range("E8").Value = "  -0.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -5.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.05"
$ws.Range("E10").Value = "  -4.65%  "

# Row 11
$ws.Range("E11").Value = "  -2.99%  "

# Row 13
$ws.Range("E13").Value = "  -0.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.566.08"
$ws.Range("E14").Value = "  -0.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.234.59"
$ws.Range("E15").Value = "  -1.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.37"
$ws.Range("E16").Value = "  -1.32%  "

# Row 17
$ws.Range("E17").Value = "  -6.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.871.71"
$ws.Range("E18").Value = "  -0.58%  "

# Row 19
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.77"
$ws.Range("E19").Value = "  -1.05%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -5.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  -5.54%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.12"
$ws.Range("E22").Value = "  -2.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.19"
$ws.Range("E23").Value = "  -0.85%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -4.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.30%  "

# Row 26
$ws.Range("E26").Value = "  -6.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.28"
$ws.Range("E27").Value = "  +2.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  -4.19%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "151.60"
$ws.Range("E30").Value = "  -1.03%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.19"
$ws.Range("E31").Value = "  -3.80%  "

# Row 32
$ws.Range("E32").Value = "  -9.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  -5.72%  "

# Row 35
$ws.Range("E35").Value = "  -1.61%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("E36").Value = "  -6.32%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.84"
$ws.Range("E37").Value = "  -7.89%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.68"
$ws.Range("E38").Value = "  -6.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0302"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  -6.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("E41").Value = "  -4.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.37"
$ws.Range("E42").Value = "  -10.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.93%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.804.35"
$ws.Range("E44").Value = "  +0.58%  "

# Row 45
$ws.Range("E45").Value = "  +13.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.184"
$ws.Range("E46").Value = "  -3.60%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.86"
$ws.Range("E47").Value = "  -3.37%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.39"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.81"
$ws.Range("E49").Value = "  -3.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.78"
$ws.Range("E50").Value = "  -7.34%  "

# Row 51
$ws.Range("E51").Value = "  -5.73%  "
